$p = $ppt.ActivePresentation

# Slide 7 ("INSERTAR DATOS EN LA TABLA CLIENTE DE VENTAS") - title placeholder
$s = $p.Slides.Item(7)
$titleShape = $s.Shapes.Item(1)
$tr = $titleShape.TextFrame.TextRange

# Split the single title run into two runs:
#   "INSERTAR DATOS EN LA TABLA CLIENTE " + "DE PEDIDOS"
$firstRun = $tr.Runs(1)
$firstRun.Text = "INSERTAR DATOS EN LA TABLA CLIENTE "
$secondRun = $tr.InsertAfter("DE PEDIDOS")
$secondRun.Font.Bold = $true
